# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 21:15:00"
$wsZhCn.Range("H2").Value = "2016-03-22 21:15:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 21:15:06"
$wsDeDe.Range("H2").Value = "2016-03-22 21:15:40"
